$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Serping1"
$ws.Cells.Item(2, 3).Value = "Sele"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 6.159891999999999
$ws.Cells.Item(2, 8).Value = 18.479676
$ws.Cells.Item(2, 9).Value = 0.007079533182016282
$ws.Cells.Item(2, 10).Value = 0.007079533182016282
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 6.490547666666667
$ws.Cells.Item(2, 14).Value = 19.471643
$ws.Cells.Item(2, 15).Value = 0.8021666724616637
$ws.Cells.Item(2, 16).Value = 0.8021666724616636
$ws.Cells.Item(2, 17).Value = 39.98107264751867
$ws.Cells.Item(2, 18).Value = 359.829653827668
$ws.Cells.Item(2, 19).Value = 0.005678965575199935
$ws.Cells.Item(2, 20).Value = 0.005678965575199934

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Serping1"
$ws.Cells.Item(3, 3).Value = "Sele"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 6.159891999999999
$ws.Cells.Item(3, 8).Value = 18.479676
$ws.Cells.Item(3, 9).Value = 0.007079533182016282
$ws.Cells.Item(3, 10).Value = 0.007079533182016282
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 1.600723
$ws.Cells.Item(3, 14).Value = 4.802169
$ws.Cells.Item(3, 15).Value = 0.1978333275383364
$ws.Cells.Item(3, 16).Value = 0.1978333275383364
$ws.Cells.Item(3, 17).Value = 9.860280801916
$ws.Cells.Item(3, 18).Value = 88.74252721724399
$ws.Cells.Item(3, 19).Value = 0.001400567606816348
$ws.Cells.Item(3, 20).Value = 0.001400567606816348

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Serping1"
$ws.Cells.Item(4, 3).Value = "Sele"
$ws.Cells.Item(4, 4).Value = "ECs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 771.4717003333334
$ws.Cells.Item(4, 8).Value = 2314.415101
$ws.Cells.Item(4, 9).Value = 0.8866485810946614
$ws.Cells.Item(4, 10).Value = 0.8866485810946614
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 6.490547666666667
$ws.Cells.Item(4, 14).Value = 19.471643
$ws.Cells.Item(4, 15).Value = 0.8021666724616637
$ws.Cells.Item(4, 16).Value = 0.8021666724616636
$ws.Cells.Item(4, 17).Value = 5007.273844497883
$ws.Cells.Item(4, 18).Value = 45065.46460048095
$ws.Cells.Item(4, 19).Value = 0.7112399419395601
$ws.Cells.Item(4, 20).Value = 0.7112399419395601

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Serping1"
$ws.Cells.Item(5, 3).Value = "Sele"
$ws.Cells.Item(5, 4).Value = "FAPs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 771.4717003333334
$ws.Cells.Item(5, 8).Value = 2314.415101
$ws.Cells.Item(5, 9).Value = 0.8866485810946614
$ws.Cells.Item(5, 10).Value = 0.8866485810946614
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 1.600723
$ws.Cells.Item(5, 14).Value = 4.802169
$ws.Cells.Item(5, 15).Value = 0.1978333275383364
$ws.Cells.Item(5, 16).Value = 0.1978333275383364
$ws.Cells.Item(5, 17).Value = 1234.912494572674
$ws.Cells.Item(5, 18).Value = 11114.21245115407
$ws.Cells.Item(5, 19).Value = 0.1754086391551014
$ws.Cells.Item(5, 20).Value = 0.1754086391551014

# Row 6
$ws.Cells.Item(6, 1).Value = "M2"
$ws.Cells.Item(6, 2).Value = "Serping1"
$ws.Cells.Item(6, 3).Value = "Sele"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.4146736666666667
$ws.Cells.Item(6, 8).Value = 1.244021
$ws.Cells.Item(6, 9).Value = 0.0004765823788590816
$ws.Cells.Item(6, 10).Value = 0.0004765823788590817
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 6.490547666666667
$ws.Cells.Item(6, 14).Value = 19.471643
$ws.Cells.Item(6, 15).Value = 0.8021666724616637
$ws.Cells.Item(6, 16).Value = 0.8021666724616636
$ws.Cells.Item(6, 17).Value = 2.691459199611445
$ws.Cells.Item(6, 18).Value = 24.223132796503
$ws.Cells.Item(6, 19).Value = 0.0003822985010032535
$ws.Cells.Item(6, 20).Value = 0.0003822985010032535

# Row 7
$ws.Cells.Item(7, 1).Value = "M2"
$ws.Cells.Item(7, 2).Value = "Serping1"
$ws.Cells.Item(7, 3).Value = "Sele"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.4146736666666667
$ws.Cells.Item(7, 8).Value = 1.244021
$ws.Cells.Item(7, 9).Value = 0.0004765823788590816
$ws.Cells.Item(7, 10).Value = 0.0004765823788590817
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 1.600723
$ws.Cells.Item(7, 14).Value = 4.802169
$ws.Cells.Item(7, 15).Value = 0.1978333275383364
$ws.Cells.Item(7, 16).Value = 0.1978333275383364
$ws.Cells.Item(7, 17).Value = 0.6637776757276667
$ws.Cells.Item(7, 18).Value = 5.973999081549001
$ws.Cells.Item(7, 19).Value = 0.00009428387785582822
$ws.Cells.Item(7, 20).Value = 0.00009428387785582822

# Row 8
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Serping1"
$ws.Cells.Item(8, 3).Value = "Sele"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 92.05234666666666
$ws.Cells.Item(8, 8).Value = 276.15704
$ws.Cells.Item(8, 9).Value = 0.1057953033444633
$ws.Cells.Item(8, 10).Value = 0.1057953033444633
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 6.490547666666667
$ws.Cells.Item(8, 14).Value = 19.471643
$ws.Cells.Item(8, 15).Value = 0.8021666724616637
$ws.Cells.Item(8, 16).Value = 0.8021666724616636
$ws.Cells.Item(8, 17).Value = 597.4701438685245
$ws.Cells.Item(8, 18).Value = 5377.23129481672
$ws.Cells.Item(8, 19).Value = 0.08486546644590044
$ws.Cells.Item(8, 20).Value = 0.08486546644590043

# Row 9
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Serping1"
$ws.Cells.Item(9, 3).Value = "Sele"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 92.05234666666666
$ws.Cells.Item(9, 8).Value = 276.15704
$ws.Cells.Item(9, 9).Value = 0.1057953033444633
$ws.Cells.Item(9, 10).Value = 0.1057953033444633
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 1.600723
$ws.Cells.Item(9, 14).Value = 4.802169
$ws.Cells.Item(9, 15).Value = 0.1978333275383364
$ws.Cells.Item(9, 16).Value = 0.1978333275383364
$ws.Cells.Item(9, 17).Value = 147.3503085133067
$ws.Cells.Item(9, 18).Value = 1326.15277661976
$ws.Cells.Item(9, 19).Value = 0.02092983689856286
$ws.Cells.Item(9, 20).Value = 0.02092983689856286
